$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.508245648867
$ws.Range("B2").Value = 0.633570024496
$ws.Range("B3").Value = 0.701599288685
$ws.Range("B4").Value = 0.678846256784
$ws.Range("B5").Value = 0.701599288685
$ws.Range("B6").Value = 0.633570024496
$ws.Range("B7").Value = 0.508245648867
$ws.Range("B8").Value = 0.545710923476
$ws.Range("B9").Value = 0.7539584878659999
$ws.Range("B10").Value = 1.16332001126
$ws.Range("B11").Value = 0.822605739677
$ws.Range("B12").Value = 1.21126960749
$ws.Range("B13").Value = 0.822605739677
$ws.Range("B14").Value = 1.16332001126
$ws.Range("B15").Value = 0.7539584878659999
$ws.Range("B16").Value = 0.545710923476
$ws.Range("B17").Value = 0.545710923476
$ws.Range("B18").Value = 0.687179174242
$ws.Range("B19").Value = 0.984686311592
$ws.Range("B20").Value = 1.23228780065
$ws.Range("B21").Value = 1.00188841038
$ws.Range("B22").Value = 1.23228780065
$ws.Range("B23").Value = 0.984686311592
$ws.Range("B24").Value = 0.687179174242
$ws.Range("B25").Value = 0.545710923476
$ws.Range("B26").Value = 0.7539584878659999
$ws.Range("B27").Value = 1.29136388046
$ws.Range("B28").Value = 1.09169830677
$ws.Range("B29").Value = 1.31079628805
$ws.Range("B30").Value = 1.09169830677
$ws.Range("B31").Value = 1.29136388046
$ws.Range("B32").Value = 0.7539584878659999
$ws.Range("B33").Value = 0.508245648867
$ws.Range("B34").Value = 1.33211027851
$ws.Range("B35").Value = 1.11029862242
$ws.Range("B36").Value = 1.34781485183
$ws.Range("B37").Value = 1.12542156452
$ws.Range("B38").Value = 1.34781485183
$ws.Range("B39").Value = 1.11029862242
$ws.Range("B40").Value = 1.33211027851
$ws.Range("B41").Value = 0.508245648867
$ws.Range("B42").Value = 0.633570024496
$ws.Range("B43").Value = 1.16332001126
$ws.Range("B44").Value = 0.984686311592
$ws.Range("B45").Value = 1.29136388046
$ws.Range("B46").Value = 1.11029862242
$ws.Range("B47").Value = 1.3798652656
$ws.Range("B48").Value = 1.11355377866
$ws.Range("B49").Value = 1.33892010145
$ws.Range("B50").Value = 1.11355377866
$ws.Range("B51").Value = 1.3798652656
$ws.Range("B52").Value = 1.11029862242
$ws.Range("B53").Value = 1.29136388046
$ws.Range("B54").Value = 0.984686311592
$ws.Range("B55").Value = 1.16332001126
$ws.Range("B56").Value = 0.633570024496
$ws.Range("B57").Value = 0.701599288685
$ws.Range("B58").Value = 0.822605739677
$ws.Range("B59").Value = 1.23228780065
$ws.Range("B60").Value = 1.09169830677
$ws.Range("B61").Value = 1.34781485183
$ws.Range("B62").Value = 1.11355377866
$ws.Range("B63").Value = 1.30139878445
$ws.Range("B64").Value = 1.03087691513
$ws.Range("B65").Value = 1.30139878445
$ws.Range("B66").Value = 1.11355377866
$ws.Range("B67").Value = 1.34781485183
$ws.Range("B68").Value = 1.09169830677
$ws.Range("B69").Value = 1.23228780065
$ws.Range("B70").Value = 0.822605739677
$ws.Range("B71").Value = 0.701599288685
$ws.Range("B72").Value = 0.678846256784
$ws.Range("B73").Value = 1.21126960749
$ws.Range("B74").Value = 1.00188841038
$ws.Range("B75").Value = 1.31079628805
$ws.Range("B76").Value = 1.12542156452
$ws.Range("B77").Value = 1.33892010145
$ws.Range("B78").Value = 1.03087691513
$ws.Range("B79").Value = 1.03087691513
$ws.Range("B80").Value = 1.33892010145
$ws.Range("B81").Value = 1.12542156452
$ws.Range("B82").Value = 1.31079628805
$ws.Range("B83").Value = 1.00188841038
$ws.Range("B84").Value = 1.21126960749
$ws.Range("B85").Value = 0.678846256784
$ws.Range("B86").Value = 0.701599288685
$ws.Range("B87").Value = 0.822605739677
$ws.Range("B88").Value = 1.23228780065
$ws.Range("B89").Value = 1.09169830677
$ws.Range("B90").Value = 1.34781485183
$ws.Range("B91").Value = 1.11355377866
$ws.Range("B92").Value = 1.30139878445
$ws.Range("B93").Value = 1.03087691513
$ws.Range("B94").Value = 1.30139878445
$ws.Range("B95").Value = 1.11355377866
$ws.Range("B96").Value = 1.34781485183
$ws.Range("B97").Value = 1.09169830677
$ws.Range("B98").Value = 1.23228780065
$ws.Range("B99").Value = 0.822605739677
$ws.Range("B100").Value = 0.701599288685
$ws.Range("B101").Value = 0.633570024496
$ws.Range("B102").Value = 1.16332001126
$ws.Range("B103").Value = 0.984686311592
$ws.Range("B104").Value = 1.29136388046
$ws.Range("B105").Value = 1.11029862242
$ws.Range("B106").Value = 1.3798652656
$ws.Range("B107").Value = 1.11355377866
$ws.Range("B108").Value = 1.33892010145
$ws.Range("B109").Value = 1.11355377866
$ws.Range("B110").Value = 1.3798652656
$ws.Range("B111").Value = 1.11029862242
$ws.Range("B112").Value = 1.29136388046
$ws.Range("B113").Value = 0.984686311592
$ws.Range("B114").Value = 1.16332001126
$ws.Range("B115").Value = 0.633570024496
$ws.Range("B116").Value = 0.508245648867
$ws.Range("B117").Value = 1.33211027851
$ws.Range("B118").Value = 1.11029862242
$ws.Range("B119").Value = 1.34781485183
$ws.Range("B120").Value = 1.12542156452
$ws.Range("B121").Value = 1.34781485183
$ws.Range("B122").Value = 1.11029862242
$ws.Range("B123").Value = 1.33211027851
$ws.Range("B124").Value = 0.508245648867
$ws.Range("B125").Value = 0.7539584878659999
$ws.Range("B126").Value = 1.29136388046
$ws.Range("B127").Value = 1.09169830677
$ws.Range("B128").Value = 1.31079628805
$ws.Range("B129").Value = 1.09169830677
$ws.Range("B130").Value = 1.29136388046
$ws.Range("B131").Value = 0.7539584878659999
$ws.Range("B132").Value = 0.545710923476
$ws.Range("B133").Value = 0.687179174242
$ws.Range("B134").Value = 0.984686311592
$ws.Range("B135").Value = 1.23228780065
$ws.Range("B136").Value = 1.00188841038
$ws.Range("B137").Value = 1.23228780065
$ws.Range("B138").Value = 0.984686311592
$ws.Range("B139").Value = 0.687179174242
$ws.Range("B140").Value = 0.545710923476
$ws.Range("B141").Value = 0.545710923476
$ws.Range("B142").Value = 0.7539584878659999
$ws.Range("B143").Value = 1.16332001126
$ws.Range("B144").Value = 0.822605739677
$ws.Range("B145").Value = 1.21126960749
$ws.Range("B146").Value = 0.822605739677
$ws.Range("B147").Value = 1.16332001126
$ws.Range("B148").Value = 0.7539584878659999
$ws.Range("B149").Value = 0.545710923476
$ws.Range("B150").Value = 0.508245648867
$ws.Range("B151").Value = 0.633570024496
$ws.Range("B152").Value = 0.701599288685
$ws.Range("B153").Value = 0.678846256784
$ws.Range("B154").Value = 0.701599288685
$ws.Range("B155").Value = 0.633570024496
$ws.Range("B156").Value = 0.508245648867
